$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 503
$ws1.Range("F3").Value = 5952
$ws1.Range("F6").Value = 108
$ws1.Range("F9").Value = 549

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 503
$ws4.Range("F3").Value = 5952
$ws4.Range("F7").Value = 108
$ws4.Range("F11").Value = 549
